$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 569.1
$ws.Range("J17").Value = 569.53424
$ws.Range("L17").Value = 1708.60272
$ws.Range("N17").Value = -2044.60272

$ws.Range("H74").Value = 4008.2222
$ws.Range("I74").Value = 3383.3333
$ws.Range("K74").Value = 3383.3333
$ws.Range("M74").Value = -2447.3333

$ws.Range("H77").Value = 4008.2222
$ws.Range("I77").Value = 3383.3333
$ws.Range("K77").Value = 16916.6665
$ws.Range("M77").Value = -12236.6665

$ws.Range("H80").Value = 491.33334
$ws.Range("I80").Value = 375.22223
$ws.Range("J80").Value = 665.5
$ws.Range("K80").Value = 1125.66669
$ws.Range("L80").Value = 1996.5
$ws.Range("M80").Value = -127.66669
$ws.Range("N80").Value = -3992.5

$ws.Range("H83").Value = 491.33334
$ws.Range("I83").Value = 375.22223
$ws.Range("J83").Value = 665.5
$ws.Range("K83").Value = 3377.00007
$ws.Range("L83").Value = 5989.5
$ws.Range("M83").Value = 1614.99993
$ws.Range("N83").Value = -15973.5

$ws.Range("H92").Value = 23811634
$ws.Range("I92").Value = 27779614
$ws.Range("K92").Value = 27779614
$ws.Range("M92").Value = -27778366

$ws.Range("H111").Value = 7895.5713
$ws.Range("I111").Value = 8887.4
$ws.Range("J111").Value = 5416
$ws.Range("K111").Value = 26662.2
$ws.Range("L111").Value = 16248
$ws.Range("M111").Value = -23595.2
$ws.Range("N111").Value = -22382

$ws.Range("H126").Value = 50000
$ws.Range("J126").Value = 50000
$ws.Range("L126").Value = 50000
$ws.Range("N126").Value = -59880

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null

$ws.Range("H132").Value = 4701.697
$ws.Range("I132").Value = 3513.7693
$ws.Range("K132").Value = 10541.3079
$ws.Range("M132").Value = -8011.3079

$ws.Range("H138").Value = 2760.0557
$ws.Range("I138").Value = 2939.7334
$ws.Range("J138").Value = 2631.7144
$ws.Range("K138").Value = 8819.200199999999
$ws.Range("L138").Value = 7895.1432
$ws.Range("M138").Value = -3679.200199999999
$ws.Range("N138").Value = -18175.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 514594.62
$ws.Range("I32").Value = 621162.9
$ws.Range("K32").Value = 621162.9
$ws.Range("M32").Value = -620875.9

$ws.Range("H61").Value = 4711.2
$ws.Range("I61").Value = 4582.4
$ws.Range("J61").Value = 4840
$ws.Range("K61").Value = 4582.4
$ws.Range("L61").Value = 4840
$ws.Range("M61").Value = -4370.4
$ws.Range("N61").Value = -5264

$ws.Range("H88").Value = 2601.6
$ws.Range("I88").Value = 2379.5
$ws.Range("J88").Value = 3490
$ws.Range("K88").Value = 2379.5
$ws.Range("L88").Value = 3490
$ws.Range("M88").Value = -1973.5
$ws.Range("N88").Value = -4302

$ws.Range("H91").Value = 2601.6
$ws.Range("I91").Value = 2379.5
$ws.Range("J91").Value = 3490
$ws.Range("K91").Value = 2379.5
$ws.Range("L91").Value = 3490
$ws.Range("M91").Value = -975.5
$ws.Range("N91").Value = -6298

$ws.Range("H92").Value = 78749.75
$ws.Range("J92").Value = 78749.75
$ws.Range("L92").Value = 78749.75
$ws.Range("N92").Value = -83741.75

$ws.Range("H113").Value = 36761
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null

$ws.Range("H123").Value = 30270
$ws.Range("J123").Value = 30270
$ws.Range("L123").Value = 30270
$ws.Range("N123").Value = -40070

$ws.Range("H130").Value = 19928
$ws.Range("J130").Value = 19928
$ws.Range("L130").Value = 19928
$ws.Range("N130").Value = -29968

$ws.Range("H131").Value = 39666.668
$ws.Range("J131").Value = 39666.668
$ws.Range("L131").Value = 39666.668
$ws.Range("N131").Value = -49746.668

$ws.Range("H132").Value = 3384.2458
$ws.Range("I132").Value = 2462.2144
$ws.Range("J132").Value = 5422.421
$ws.Range("K132").Value = 7386.6432
$ws.Range("L132").Value = 16267.263
$ws.Range("M132").Value = -4856.6432
$ws.Range("N132").Value = -21327.263

$ws.Range("H136").Value = 4711.2
$ws.Range("I136").Value = 4582.4
$ws.Range("J136").Value = 4840
$ws.Range("K136").Value = 13747.2
$ws.Range("L136").Value = 14520
$ws.Range("M136").Value = -11197.2
$ws.Range("N136").Value = -19620

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13735.588
$ws.Range("I82").Value = 5439.25
$ws.Range("J82").Value = 21110.111
$ws.Range("K82").Value = 5439.25
$ws.Range("L82").Value = 21110.111
$ws.Range("M82").Value = -5056.25
$ws.Range("N82").Value = -21876.111

$ws.Range("H85").Value = 13735.588
$ws.Range("I85").Value = 5439.25
$ws.Range("J85").Value = 21110.111
$ws.Range("K85").Value = 5439.25
$ws.Range("L85").Value = 21110.111
$ws.Range("M85").Value = -4113.25
$ws.Range("N85").Value = -23762.111

$ws.Range("H88").Value = 35000
$ws.Range("J88").Value = 35000
$ws.Range("L88").Value = 35000
$ws.Range("N88").Value = -35812

$ws.Range("H91").Value = 35000
$ws.Range("J91").Value = 35000
$ws.Range("L91").Value = 35000
$ws.Range("N91").Value = -37808

$ws.Range("H94").Value = 2224.111
$ws.Range("I94").Value = 2403.4
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 2403.4
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -1952.4
$ws.Range("N94").Value = -2902

$ws.Range("H105").Value = 8335572
$ws.Range("I105").Value = 12501977
$ws.Range("J105").Value = 2762.2
$ws.Range("K105").Value = 12501977
$ws.Range("L105").Value = 2762.2
$ws.Range("M105").Value = -12500230
$ws.Range("N105").Value = -6256.2

$ws.Range("H134").Value = 3504.88
$ws.Range("I134").Value = 3280.6
$ws.Range("J134").Value = 3841.3
$ws.Range("K134").Value = 9841.799999999999
$ws.Range("L134").Value = 11523.9
$ws.Range("M134").Value = -7306.799999999999
$ws.Range("N134").Value = -16593.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 388.57144
$ws.Range("I22").Value = 265
$ws.Range("J22").Value = 697.5
$ws.Range("K22").Value = 265
$ws.Range("L22").Value = 697.5
$ws.Range("M22").Value = 85
$ws.Range("N22").Value = -1397.5

$ws.Range("H96").Value = 43333.332
$ws.Range("J96").Value = 43333.332
$ws.Range("L96").Value = 43333.332
$ws.Range("N96").Value = -48825.332

$ws.Range("H107").Value = 2976785
$ws.Range("I107").Value = 5682294
$ws.Range("J107").Value = 725
$ws.Range("K107").Value = 5682294
$ws.Range("L107").Value = 725
$ws.Range("M107").Value = -5680374
$ws.Range("N107").Value = -4565

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 673.1667
$ws.Range("I13").Value = 150
$ws.Range("J13").Value = 934.75
$ws.Range("K13").Value = 450
$ws.Range("L13").Value = 2804.25
$ws.Range("M13").Value = -282
$ws.Range("N13").Value = -3140.25

$ws.Range("H113").Value = 1151.5714
$ws.Range("I113").Value = 657
$ws.Range("J113").Value = 1601.1818
$ws.Range("K113").Value = 1971
$ws.Range("L113").Value = 4803.5454
$ws.Range("M113").Value = 199
$ws.Range("N113").Value = -9143.545399999999

$ws.Range("H114").Value = 490.45456
$ws.Range("J114").Value = 907.25
$ws.Range("L114").Value = 2721.75
$ws.Range("N114").Value = -9229.75

$ws.Range("H129").Value = 1877.7368
$ws.Range("I129").Value = 710
$ws.Range("K129").Value = 2130
$ws.Range("M129").Value = 2870

$ws.Range("H136").Value = 3887.375
$ws.Range("I136").Value = 3833
$ws.Range("J136").Value = 3920
$ws.Range("K136").Value = 11499
$ws.Range("L136").Value = 11760
$ws.Range("M136").Value = -6399
$ws.Range("N136").Value = -21960

$ws.Range("H137").Value = 8342327.5
$ws.Range("I137").Value = 23826910
$ws.Range("J137").Value = 4475.3076
$ws.Range("K137").Value = 71480730
$ws.Range("L137").Value = 13425.9228
$ws.Range("M137").Value = -71475630
$ws.Range("N137").Value = -23625.9228

$ws.Range("H139").Value = 3389.1135
$ws.Range("J139").Value = 4210.6772
$ws.Range("L139").Value = 12632.0316
$ws.Range("N139").Value = -22912.0316

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2093.3333
$ws.Range("I97").Value = 1800
$ws.Range("J97").Value = 2350
$ws.Range("K97").Value = 1800
$ws.Range("L97").Value = 2350
$ws.Range("M97").Value = -1304
$ws.Range("N97").Value = -3342

$ws.Range("H122").Value = 1600
$ws.Range("I122").Value = 1600
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4800
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2350
$ws.Range("N122").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 35106.5
$ws.Range("I24").Value = 206
$ws.Range("K24").Value = 206
$ws.Range("M24").Value = 137

$ws.Range("H127").Value = 30846.574
$ws.Range("J127").Value = 30846.574
$ws.Range("L127").Value = 30846.574
$ws.Range("N127").Value = -40766.574

$ws.Range("H130").Value = 28271.428
$ws.Range("J130").Value = 28271.428
$ws.Range("L130").Value = 28271.428
$ws.Range("N130").Value = -38311.428

$ws.Range("H132").Value = 2444.2727
$ws.Range("I132").Value = 1922.0588
$ws.Range("J132").Value = 4219.8
$ws.Range("K132").Value = 5766.1764
$ws.Range("L132").Value = 12659.4
$ws.Range("M132").Value = -3236.1764
$ws.Range("N132").Value = -17719.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4632.979
$ws.Range("J96").Value = 5137.838
$ws.Range("L96").Value = 5137.838
$ws.Range("N96").Value = -7883.838

$ws.Range("H97").Value = 90000
$ws.Range("J97").Value = 90000
$ws.Range("L97").Value = 90000
$ws.Range("N97").Value = -91982

$ws.Range("H123").Value = 22700
$ws.Range("J123").Value = 22700
$ws.Range("L123").Value = 22700
$ws.Range("N123").Value = -32500

$ws.Range("H128").Value = 48500
$ws.Range("J128").Value = 48500
$ws.Range("L128").Value = 48500
$ws.Range("N128").Value = -58460

$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

$ws.Range("H132").Value = 5378004.5
$ws.Range("I132").Value = 1236.0834
$ws.Range("K132").Value = 3708.2502
$ws.Range("M132").Value = -1178.2502
